$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:132 down to 60:133
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new data record
$ws.Range("A59").Value = 6
$ws.Range("B59").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = 44467
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = 100112022
$ws.Range("G59").Value = "Arveja Verde"
$ws.Range("H59").Value = "Perfection"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 300
$ws.Range("K59").Value = 24000
$ws.Range("L59").Value = 25000
$ws.Range("M59").Value = 24400
$ws.Range("N59").Value = "`$/malla 25 kilos"
$ws.Range("O59").Value = "Provincia de Huasco"
$ws.Range("P59").Value = 976
$ws.Range("Q59").Value = 25
$ws.Range("R59").Value = "Hortaliza"
